$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 314 (pushes existing rows 314:433 down to 315:434)
$ws.Rows.Item(314).Insert()

# Populate the newly inserted row 314 with the new weekly record
$ws.Range("A314").Value = 4
$ws.Range("B314").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C314").Value = "Los Lagos"
$ws.Range("D314").Value = 45146
$ws.Range("E314").Value = 10
$ws.Range("F314").Value = 100112044
$ws.Range("G314").Value = "Perejil"
$ws.Range("H314").Value = "Sin especificar"
$ws.Range("I314").Value = "Primera"
$ws.Range("J314").Value = 180
$ws.Range("K314").Value = 6000
$ws.Range("L314").Value = 6000
$ws.Range("M314").Value = 6000
$ws.Range("N314").Value = "`$/docena de atados (3 kilos)"
$ws.Range("O314").Value = "Región Metropolitana"
$ws.Range("P314").Value = 2000
$ws.Range("Q314").Value = 3
$ws.Range("R314").Value = "Hortaliza"
